$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Punto 3 (row 4) - list of ints, stored as text
$ws.Range("B4").Value = "[63, 12, 6, 20, 8]"

# Punto 5 - Máximo (row 7) - numeric
$ws.Range("B7").Value = 0.9846832506704628

# Punto 5 - Mínimo (row 8) - numeric
$ws.Range("B8").Value = 0.0008586509496919525

# Punto 5 - Índice Máx (row 9) - numeric
$ws.Range("B9").Value = 77

# Punto 5 - Índice Mín (row 10) - numeric
$ws.Range("B10").Value = 15

# Punto 7 (row 12) - matrix, stored as text
$ws.Range("B12").Value = "[[79, 34], [52, 65]]"

# Punto 9 (row 14) - matrix, stored as text
$ws.Range("B14").Value = "[[64, 83, 70], [95, 26, 47], [74, 59, 97]]"

# Punto 10 (row 15) - list of floats, stored as text
$ws.Range("B15").Value = "[0.669980916476668, 0.900369825570438, 0.9729553122022809, 0.5357689696200466, 0.7412527302268486]"
